$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 693.7692
$ws.Range("I33").Value = 1330.5
$ws.Range("J33").Value = 148
$ws.Range("K33").Value = 1330.5
$ws.Range("L33").Value = 148
$ws.Range("M33").Value = -1101.5
$ws.Range("N33").Value = -606

$ws.Range("H40").Value = 1360
$ws.Range("I40").Value = 1450
$ws.Range("J40").Value = 1300
$ws.Range("K40").Value = 1450
$ws.Range("L40").Value = 1300
$ws.Range("M40").Value = -1275
$ws.Range("N40").Value = -1650

$ws.Range("H100").Value = 4237.7144
$ws.Range("I100").Value = 2201.25
$ws.Range("J100").Value = 6953
$ws.Range("K100").Value = 2201.25
$ws.Range("L100").Value = 6953
$ws.Range("M100").Value = -1660.25
$ws.Range("N100").Value = -8035

$ws.Range("H138").Value = 3690.7815
$ws.Range("I138").Value = 3077.0435
$ws.Range("J138").Value = 3911.3438
$ws.Range("K138").Value = 9231.130500000001
$ws.Range("L138").Value = 11734.0314
$ws.Range("M138").Value = -4091.130500000001
$ws.Range("N138").Value = -22014.0314

$ws.Range("H140").Value = 72993.19
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 72993.19
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 72993.19
$ws.Range("N140").Value = -83353.19

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5860.28
$ws.Range("I32").Value = 4348.1304
$ws.Range("J32").Value = 23250
$ws.Range("K32").Value = 4348.1304
$ws.Range("L32").Value = 23250
$ws.Range("M32").Value = -4061.1304
$ws.Range("N32").Value = -23824

$ws.Range("H61").Value = 10756119
$ws.Range("I61").Value = 19609906
$ws.Range("J61").Value = 5092.857
$ws.Range("K61").Value = 19609906
$ws.Range("L61").Value = 5092.857
$ws.Range("M61").Value = -19609694
$ws.Range("N61").Value = -5516.857

$ws.Range("H122").Value = 78772.38
$ws.Range("I122").Value = 92558.27
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 277674.81
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -275224.81
$ws.Range("N122").Value = -13750

$ws.Range("H136").Value = 10756119
$ws.Range("I136").Value = 19609906
$ws.Range("J136").Value = 5092.857
$ws.Range("K136").Value = 58829718
$ws.Range("L136").Value = 15278.571
$ws.Range("M136").Value = -58827168
$ws.Range("N136").Value = -20378.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1774.1
$ws.Range("I99").Value = 1461.4286
$ws.Range("J99").Value = 2503.6667
$ws.Range("K99").Value = 1461.4286
$ws.Range("L99").Value = 2503.6667
$ws.Range("M99").Value = 36.57140000000004
$ws.Range("N99").Value = -5499.6667

$ws.Range("H134").Value = 2448.186
$ws.Range("I134").Value = 2339.4055
$ws.Range("J134").Value = 3119
$ws.Range("K134").Value = 7018.2165
$ws.Range("L134").Value = 9357
$ws.Range("M134").Value = -4483.2165
$ws.Range("N134").Value = -14427

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2451.6897
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 2503.8076
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 2503.8076
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -5499.8076

$ws.Range("H126").Value = 2451.6897
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2503.8076
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 7511.4228
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -12451.4228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 174.66667
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 209.09091
$ws.Range("K2").Value = 480
$ws.Range("L2").Value = 1254.54546
$ws.Range("M2").Value = -367
$ws.Range("N2").Value = -1480.54546

$ws.Range("H26").Value = 5399.775
$ws.Range("I26").Value = 69.3
$ws.Range("J26").Value = 7176.6
$ws.Range("K26").Value = 207.9
$ws.Range("L26").Value = 21529.8
$ws.Range("M26").Value = 80.10000000000002
$ws.Range("N26").Value = -22105.8

$ws.Range("H63").Value = 8000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 0
$ws.Range("N63").Value = -25498
$ws.Range("L63").ClearContents()
$ws.Range("M63").ClearContents()

$ws.Range("H64").Value = 2915.7856
$ws.Range("I64").Value = 1150
$ws.Range("J64").Value = 3051.6155
$ws.Range("K64").Value = 3450
$ws.Range("L64").Value = 9154.8465
$ws.Range("M64").Value = -3180
$ws.Range("N64").Value = -9694.8465

$ws.Range("H66").Value = 8000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 0
$ws.Range("N66").Value = -79488
$ws.Range("L66").ClearContents()
$ws.Range("M66").ClearContents()

$ws.Range("H67").Value = 2915.7856
$ws.Range("I67").Value = 1150
$ws.Range("J67").Value = 3051.6155
$ws.Range("K67").Value = 3450
$ws.Range("L67").Value = 9154.8465
$ws.Range("M67").Value = -2514
$ws.Range("N67").Value = -11026.8465

$ws.Range("H138").Value = 2800.3044
$ws.Range("I138").Value = 1448.2632
$ws.Range("J138").Value = 9222.5
$ws.Range("K138").Value = 4344.7896
$ws.Range("L138").Value = 27667.5
$ws.Range("M138").Value = 795.2103999999999
$ws.Range("N138").Value = -37947.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 70546.56
$ws.Range("I113").Value = 86233.234
$ws.Range("J113").Value = 2571
$ws.Range("K113").Value = 86233.234
$ws.Range("L113").Value = 2571
$ws.Range("M113").Value = -84063.234
$ws.Range("N113").Value = -6911

$ws.Range("H123").Value = 8589.286
$ws.Range("I123").Value = 3000
$ws.Range("J123").Value = 9904.412
$ws.Range("K123").Value = 3000
$ws.Range("L123").Value = 9904.412
$ws.Range("M123").Value = -550
$ws.Range("N123").Value = -14804.412

$ws.Range("H126").Value = 5120
$ws.Range("I126").Value = 4660
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 13980
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -11510
$ws.Range("N126").Value = -24440

$ws.Range("H132").Value = 55564292
$ws.Range("I132").Value = 83343940
$ws.Range("J132").Value = 5004.1665
$ws.Range("K132").Value = 250031820
$ws.Range("L132").Value = 15012.4995
$ws.Range("M132").Value = -250029290
$ws.Range("N132").Value = -20072.4995

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").ClearContents()

$ws.Range("H141").Value = 73146.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 73146.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 73146.5
$ws.Range("N141").Value = -83506.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9585.625
$ws.Range("I7").Value = 9120
$ws.Range("J7").Value = 10051.25
$ws.Range("K7").Value = 9120
$ws.Range("L7").Value = 10051.25
$ws.Range("M7").Value = -9008
$ws.Range("N7").Value = -10275.25

$ws.Range("H16").Value = 247.11111
$ws.Range("I16").Value = 247.11111
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 247.11111
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

$ws.Range("H40").Value = 3347.6365
$ws.Range("I40").Value = 2979.889
$ws.Range("J40").Value = 5002.5
$ws.Range("K40").Value = 2979.889
$ws.Range("L40").Value = 5002.5
$ws.Range("M40").Value = -2843.889
$ws.Range("N40").Value = -5274.5

$ws.Range("H96").Value = 20000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 20000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 20000
$ws.Range("N96").Value = -25492

$ws.Range("H100").Value = 73003.38
$ws.Range("I100").Value = 131845.72
$ws.Range("J100").Value = 4354
$ws.Range("K100").Value = 131845.72
$ws.Range("L100").Value = 4354
$ws.Range("M100").Value = -131304.72
$ws.Range("N100").Value = -5436

$ws.Range("H126").Value = 9585.625
$ws.Range("I126").Value = 9120
$ws.Range("J126").Value = 10051.25
$ws.Range("K126").Value = 27360
$ws.Range("L126").Value = 30153.75
$ws.Range("M126").Value = -24890
$ws.Range("N126").Value = -35093.75

$ws.Range("H140").Value = 65771.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 65771.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 65771.5
$ws.Range("N140").Value = -76131.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 39490
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 39490
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 39490
$ws.Range("N99").Value = -45480

$ws.Range("H113").Value = 1244.9333
$ws.Range("I113").Value = 1390.3077
$ws.Range("J113").Value = 300
$ws.Range("K113").Value = 4170.9231
$ws.Range("L113").Value = 900
$ws.Range("M113").Value = -2000.9231
$ws.Range("N113").Value = -5240

$ws.Range("H122").Value = 2873.1538
$ws.Range("I122").Value = 2053.4666
$ws.Range("J122").Value = 3990.9092
$ws.Range("K122").Value = 6160.399800000001
$ws.Range("L122").Value = 11972.7276
$ws.Range("M122").Value = -3710.399800000001
$ws.Range("N122").Value = -16872.7276

$ws.Range("H126").Value = 1861.875
$ws.Range("I126").Value = 2100.0908
$ws.Range("J126").Value = 1337.8
$ws.Range("K126").Value = 6300.2724
$ws.Range("L126").Value = 4013.4
$ws.Range("M126").Value = -3830.2724
$ws.Range("N126").Value = -8953.4

$ws.Range("H136").Value = 4666.212
$ws.Range("I136").Value = 5938.5
$ws.Range("J136").Value = 3939.1904
$ws.Range("K136").Value = 17815.5
$ws.Range("L136").Value = 11817.5712
$ws.Range("M136").Value = -15265.5
$ws.Range("N136").Value = -16917.5712

$ws.Range("H141").Value = 66610.81
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 66610.81
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 66610.81
$ws.Range("N141").Value = -76970.81
